# Aplicación - Actividad terminada. TODO: Validaciones
#
# Update the single data row of Carga_Tarifas:
#  - Mes (A2):      1    -> 2
#  - Año (B2):      2020 -> 2020 (unchanged)
#  - Tipo (C2):     "I"  -> "R"
#  - Tarifa_B (D2): 12   -> 1.55
#  - Tarifa_I (E2): 12   -> 2
#  - Tarifa_E (F2): 45   -> 2.81
# and switch the row's alignment from centered to left-aligned, then move
# the active selection to H5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRow = $ws.Range("A2:F2")
$dataRow.HorizontalAlignment = -4131  # xlLeft

$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 2020
$ws.Range("C2").Value = "R"
$ws.Range("D2").Value = 1.55
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 2.81

$ws.Range("H5").Select() | Out-Null
